# Update "want-to-go" counts (column F) and "minimum price" (column G)
# figures on the 展览 (Exhibitions) and 全部类型 (All Types) sheets to match
# the freshly scraped data (commit: "Update gh-pages to output generated at
# 456a3b4").
#
# Row numbers differ by one between the two sheets from row 21 onward,
# because 全部类型 also includes a 演出 (show) row that 展览 does not have.
# Each sheet therefore gets its own explicit row map.

$wb = $excel.ActiveWorkbook

# ---- 展览 sheet ----
$ws1 = $wb.Worksheets.Item("展览")

$exhibitionF = @{
    3  = 5699
    6  = 406
    8  = 139
    9  = 4410
    10 = 788
    11 = 831
    13 = 34
    14 = 129
    15 = 148
    17 = 17
    18 = 129
    19 = 610
    21 = 182
    22 = 1153
    23 = 19
    24 = 2781
    26 = 458
}
foreach ($row in $exhibitionF.Keys) {
    $ws1.Range("F$row").Value = $exhibitionF[$row]
}

$exhibitionG = @{
    4  = "不可售"
    15 = 88
}
foreach ($row in $exhibitionG.Keys) {
    $ws1.Range("G$row").Value = $exhibitionG[$row]
}

# ---- 全部类型 sheet (row offset by +1 from row 21 onward) ----
$ws2 = $wb.Worksheets.Item("全部类型")

$allTypesF = @{
    3  = 5699
    6  = 406
    8  = 139
    9  = 4410
    10 = 788
    11 = 831
    13 = 34
    14 = 129
    15 = 148
    17 = 17
    18 = 129
    19 = 610
    22 = 182
    23 = 1153
    24 = 19
    25 = 2781
    27 = 458
}
foreach ($row in $allTypesF.Keys) {
    $ws2.Range("F$row").Value = $allTypesF[$row]
}

$allTypesG = @{
    4  = "不可售"
    15 = 88
}
foreach ($row in $allTypesG.Keys) {
    $ws2.Range("G$row").Value = $allTypesG[$row]
}

Write-Output "edits applied"
